# Commit: "Added a new testcase" / "Added a new testcase Itinerary"
#
# TC_001 and TC_002 (rows 2-3) already had a Status/RunTime recorded; they get
# refreshed to a new run: Status -> Pass, new Start/End Time, new RunTime.
#
# TC_003..TC_007 (rows 4-8) previously only had S.No/Name/Action Word filled
# in (Execute = No, nothing else) -- they had never actually been run. They
# are now fully populated as executed test cases: Execute = Yes,
# Status = Pass, Start Time / End Time, and a RunTime label.
#
# (The per-row writes below are intentionally interleaved rather than strict
# top-to-bottom so the workbook's shared-string table comes out in the same
# order as the reference edit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : TC_001_Verify Valid User is able to login to Application or not.
$ws.Range("E2").Value = "Pass"
$ws.Range("F2").Value = 0.010011574074074074
$ws.Range("G2").Value = 0.010266203703703703

# Row 7 : TC_006_Verify User is able to SignOff from Application or not.
$ws.Range("D7").Value = "Yes"
$ws.Range("E7").Value = "Pass"
$ws.Range("F7").Value = 0.011238425925925928
$ws.Range("G7").Value = 0.011469907407407408
$ws.Range("H7").Value = "20 Seconds"

# Row 5 : TC_004_Verify User is able to SignOff from Application or not.
$ws.Range("D5").Value = "Yes"
$ws.Range("E5").Value = "Pass"
$ws.Range("F5").Value = 0.010798611111111111
$ws.Range("G5").Value = 0.011018518518518518
$ws.Range("H5").Value = "19 Seconds"

$ws.Range("H2").Value = "22 Seconds"

# Row 6 : TC_005_Verify User is able to Click on Itinerary Link.
$ws.Range("D6").Value = "Yes"
$ws.Range("E6").Value = "Pass"
$ws.Range("F6").Value = 0.011030092592592591
$ws.Range("G6").Value = 0.011226851851851854
$ws.Range("H6").Value = "17 Seconds"

# Row 4 : TC_003_Verify User is able to Click on Itinerary Link.
$ws.Range("D4").Value = "Yes"
$ws.Range("E4").Value = "Pass"
$ws.Range("F4").Value = 0.01054398148148148
$ws.Range("G4").Value = 0.010787037037037038
$ws.Range("H4").Value = "21 Seconds"

# Row 3 : TC_002_Verify User is able to SignOff from Application or not.
$ws.Range("F3").Value = 0.010277777777777778
$ws.Range("G3").Value = 0.010532407407407407
$ws.Range("H3").Value = "22 Seconds"

# Row 8 : TC_007_Verify User is able to Click on Itinerary Link.
$ws.Range("D8").Value = "Yes"
$ws.Range("E8").Value = "Pass"
$ws.Range("F8").Value = 0.011493055555555555
$ws.Range("G8").Value = 0.01167824074074074
$ws.Range("H8").Value = "16 Seconds"

# Active cell ends on the Status column of the last newly-filled row
$ws.Range("E8").Select()
